# Generate Report for Handoff
# - Overview!G4:G7       : bump "Latest HO Xliff Generate Date" to the new run timestamp
# - zh-cn!E4:E7 (Priority): "low" -> "ht" for the rows that are "Ready for handoff"
# - zh-cn!H4:H7           : bump "Latest Handoff Datetime" for the same rows
# - de-de!E4:E7 (Priority): "low" -> "ht" for the rows that are "Ready for handoff"
# - de-de!H4:H7           : bump "Latest Handoff Datetime" for the same rows

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4:G7").Value = "2016-08-17 00:29:46"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-17 00:29:41"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-17 00:29:46"
